# Applies the assignment3-1.docx edit:
#  - wraps a few technical/camelCase tokens in spell-check proofErr markers
#    (ie, sendEmail, getEmail, EmailGenerationSystem, EmaiLGenerationSystem,
#    boolean)
#  - appends two new sentences (singleton note + encryption note)
#
# Strategy: Range.InsertXML (pkg:package payload) replaces a range's
# underlying OOXML wholesale, which is what lets us split plain runs and
# splice in <w:proofErr/> spellStart/spellEnd markers the way Word's own
# spell-checker would -- something plain Find/Replace text substitution
# cannot do. In this runtime InsertXML only behaves as an exact in-place
# swap when the target Range spans a *whole paragraph* (start of its text
# through its trailing paragraph mark); collapsing/clearing a sub-run range
# first and then calling InsertXML on it ends up replacing the whole
# enclosing paragraph anyway, so we target whole paragraphs directly and
# feed back a complete, edited copy of that paragraph's own original OOXML
# (attributes like w14:paraId / w:rsidR are preserved because we start from
# the exact original markup and only substitute the specific run(s) that
# change).

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Get-ParagraphByText($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    throw "Paragraph not found containing: $needle"
}

function Set-ParagraphXml($paragraph, $newParagraphXml) {
    $paragraph.Range.InsertXML($pkgOpen + '<w:body>' + $newParagraphXml + '</w:body>' + $pkgClose)
}

# --- Paragraph: "In terms of flexibility, ... in this way." --------------
$p1 = Get-ParagraphByText("in subclasses that represent these types")
$p1Xml = '<w:p w14:paraId="49C6EEEA" w14:textId="79A88438" w:rsidR="00F724A9" w:rsidRDefault="00F724A9" w:rsidP="00F724A9"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>In terms of flexibility, it is my goal to ensure that this program is as flexible as possible. Therefore, the program will</w:t></w:r><w:r w:rsidR="00A2153E"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>utilize a factory design pattern. The factory will be for the customers. All of the common code among all of the customers</w:t></w:r><w:r w:rsidR="00A2153E"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>will be placed in a separate Customers abstract class, and the code specific to the individual types of customers will be placed</w:t></w:r><w:r w:rsidR="00A2153E"><w:t xml:space="preserve"> </w:t></w:r>' + `
          '<w:r><w:t>in subclasses that represent these types (</w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r><w:t>ie</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r><w:t xml:space="preserve"> Newbie.java or Business.java). New customer types can easily be added in this way.</w:t></w:r>' + `
          '</w:p>'
Set-ParagraphXml $p1 $p1Xml

# --- Paragraph: "I utilized the factory design pattern ... EmailGenerationSystem." + new sentence
$p2 = Get-ParagraphByText("it will sendEmail")
$p2Xml = '<w:p w14:paraId="339D049D" w14:textId="56F8D156" w:rsidR="00F724A9" w:rsidRDefault="00F724A9" w:rsidP="00F724A9"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>I utilized the factory design pattern as that made the most logical sense. When a company creates a new message request,</w:t></w:r><w:r w:rsidR="00A2153E"><w:t xml:space="preserve"> </w:t></w:r>' + `
          '<w:r><w:t xml:space="preserve">it will </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r><w:t>sendEmail</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r><w:t xml:space="preserve">(), learn the context of the message by entering the factory and calling </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r><w:t>getEmail</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r><w:t>() on all of the subclasses</w:t></w:r>' + `
          '<w:r w:rsidR="00A2153E"><w:t xml:space="preserve"> </w:t></w:r>' + `
          '<w:r><w:t xml:space="preserve">in the factory, and populate the contents of an Email object. This will all be run by the </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r><w:t>EmailGenerationSystem</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r><w:t>.</w:t></w:r>' + `
          '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
          '<w:r><w:t>Additionally, for the system, I used a singleton so that there can only be one instance to limit confusion and duplication.</w:t></w:r>' + `
          '</w:p>'
Set-ParagraphXml $p2 $p2Xml

# --- Paragraph: "Currently, the program operates ... client types." + new sentence
$p3 = Get-ParagraphByText("EmaiLGenerationSystem runs the setup")
$p3Xml = '<w:p w14:paraId="375EDFD9" w14:textId="6DE92197" w:rsidR="004A0B27" w:rsidRDefault="00F724A9" w:rsidP="00F724A9"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Currently, the program operates from a simple factory design. A set of companies exist as well as a subset of their</w:t></w:r><w:r w:rsidR="00A2153E"><w:t xml:space="preserve"> </w:t></w:r>' + `
          '<w:r><w:t xml:space="preserve">customers, each of various types. The </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r><w:t>EmaiLGenerationSystem</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r><w:t xml:space="preserve"> runs the setup, creates a company, who sends messages based on their client types.</w:t></w:r>' + `
          '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
          '<w:r><w:t>I added a rudimentary (</w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r><w:t>boolean</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r><w:t xml:space="preserve"> based, no cryptography) encryption option for business emails.</w:t></w:r>' + `
          '</w:p>'
Set-ParagraphXml $p3 $p3Xml
